# #36225 translate object names to polish
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate English object names to Polish equivalents
$ws.Range("A2").Value = "koło"
$ws.Range("B2").Value = "kwadrat"
$ws.Range("C2").Value = "kwadrat"

$ws.Range("A3").Value = "kwadrat"
$ws.Range("B3").Value = "trójkąt"
$ws.Range("C3").Value = "trójkąt"

# Update the active selection to C3
$ws.Range("C3").Select()
